$d = $word.ActiveDocument

# The document currently ends with an empty paragraph (pStyle a7, left-indent 360,
# sz 28 / szCs 32 in its mark run properties) immediately before the section break.
# That paragraph must stay untouched; the new "payment_detail" heading, its
# description table, and a fresh trailing paragraph are appended after it.
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter() | Out-Null

# $target is the brand-new (still empty) paragraph created above; InsertXML on its
# collapsed range fills it with the first inserted block and appends the rest
# (heading paragraph, table, final paragraph) right after it.
$target = $d.Paragraphs.Last

$xml = @'

    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="a7"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:ind w:firstLineChars="0"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="32"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="32"/>
        </w:rPr>
        <w:t>p</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="32"/>
        </w:rPr>
        <w:t>ayment_detail</w:t>
      </w:r>
    </w:p>
    <w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:tblPr>
        <w:tblStyle w:val="a8"/>
        <w:tblW w:w="0" w:type="auto"/>
        <w:tblInd w:w="360" w:type="dxa"/>
        <w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/>
      </w:tblPr>
      <w:tblGrid>
        <w:gridCol w:w="1497"/>
        <w:gridCol w:w="1078"/>
        <w:gridCol w:w="1000"/>
        <w:gridCol w:w="1126"/>
        <w:gridCol w:w="3235"/>
      </w:tblGrid>
      <w:tr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="1497" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="a7"/>
              <w:ind w:firstLineChars="0" w:firstLine="0"/>
              <w:rPr>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:hint="eastAsia"/>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>Column</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="985" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="a7"/>
              <w:ind w:firstLineChars="0" w:firstLine="0"/>
              <w:rPr>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:hint="eastAsia"/>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>Type</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="1004" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="a7"/>
              <w:ind w:firstLineChars="0" w:firstLine="0"/>
              <w:rPr>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:hint="eastAsia"/>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>Default</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="1132" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="a7"/>
              <w:ind w:firstLineChars="0" w:firstLine="0"/>
              <w:rPr>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:hint="eastAsia"/>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>Nullable</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="3318" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="a7"/>
              <w:ind w:firstLineChars="0" w:firstLine="0"/>
              <w:rPr>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:hint="eastAsia"/>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>Description</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="1497" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="a7"/>
              <w:ind w:firstLineChars="0" w:firstLine="0"/>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:hint="eastAsia"/>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>u</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>id</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="985" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="a7"/>
              <w:ind w:firstLineChars="0" w:firstLine="0"/>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:hint="eastAsia"/>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>I</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>nt</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="1004" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="a7"/>
              <w:ind w:firstLineChars="0" w:firstLine="0"/>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:hint="eastAsia"/>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>/</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="1132" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="a7"/>
              <w:ind w:firstLineChars="0" w:firstLine="0"/>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:hint="eastAsia"/>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>N</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="3318" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="a7"/>
              <w:ind w:firstLineChars="0" w:firstLine="0"/>
              <w:rPr>
                <w:rFonts w:hint="eastAsia"/>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>Uid</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="1497" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="a7"/>
              <w:ind w:firstLineChars="0" w:firstLine="0"/>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>Card_number</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="985" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="a7"/>
              <w:ind w:firstLineChars="0" w:firstLine="0"/>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:hint="eastAsia"/>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>S</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>tring</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="1004" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="a7"/>
              <w:ind w:firstLineChars="0" w:firstLine="0"/>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:hint="eastAsia"/>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>/</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="1132" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="a7"/>
              <w:ind w:firstLineChars="0" w:firstLine="0"/>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:hint="eastAsia"/>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>N</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="3318" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="a7"/>
              <w:ind w:firstLineChars="0" w:firstLine="0"/>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>Card_number</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="1497" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="a7"/>
              <w:ind w:firstLineChars="0" w:firstLine="0"/>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:hint="eastAsia"/>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>A</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>ddress</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="985" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="a7"/>
              <w:ind w:firstLineChars="0" w:firstLine="0"/>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:hint="eastAsia"/>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>S</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>tring</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="1004" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="a7"/>
              <w:ind w:firstLineChars="0" w:firstLine="0"/>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:hint="eastAsia"/>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>/</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="1132" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="a7"/>
              <w:ind w:firstLineChars="0" w:firstLine="0"/>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:hint="eastAsia"/>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>N</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="3318" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="a7"/>
              <w:ind w:firstLineChars="0" w:firstLine="0"/>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>Bill address</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="1497" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="a7"/>
              <w:ind w:firstLineChars="0" w:firstLine="0"/>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>Add_time</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="985" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="a7"/>
              <w:ind w:firstLineChars="0" w:firstLine="0"/>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>Datetime</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="1004" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="a7"/>
              <w:ind w:firstLineChars="0" w:firstLine="0"/>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:hint="eastAsia"/>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>/</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="1132" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="a7"/>
              <w:ind w:firstLineChars="0" w:firstLine="0"/>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:hint="eastAsia"/>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>N</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="3318" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="a7"/>
              <w:ind w:firstLineChars="0" w:firstLine="0"/>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:hint="eastAsia"/>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>Y</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>YYYMMddHHmmss</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="1497" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="a7"/>
              <w:ind w:firstLineChars="0" w:firstLine="0"/>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:hint="eastAsia"/>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>S</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>erial</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="985" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="a7"/>
              <w:ind w:firstLineChars="0" w:firstLine="0"/>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="1004" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="a7"/>
              <w:ind w:firstLineChars="0" w:firstLine="0"/>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="1132" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="a7"/>
              <w:ind w:firstLineChars="0" w:firstLine="0"/>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="3318" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="a7"/>
              <w:ind w:firstLineChars="0" w:firstLine="0"/>
              <w:rPr>
              </w:rPr>
            </w:pPr>
          </w:p>
        </w:tc>
      </w:tr>
    </w:tbl>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="a7"/>
        <w:ind w:left="360" w:firstLineChars="0" w:firstLine="0"/>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="32"/>
        </w:rPr>
      </w:pPr>
    </w:p>

'@

$target.Range.InsertXML($xml) | Out-Null

Write-Output "Inserted payment_detail heading, table and trailing paragraph"
